$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-02-22 18:00:43"
$wsZhCn.Range("G5").Value = "2016-02-22 18:01:57"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-02-22 18:01:02"
$wsDeDe.Range("G5").Value = "2016-02-22 18:02:18"
